$wb = $excel.ActiveWorkbook

# Remember which sheet was originally active so it can be restored at the end.
$originalActiveSheet = $wb.ActiveSheet

# --- NewLoanInput sheet: move selection from B2 to B10 ---
$wsLoanInput = $wb.Worksheets.Item("NewLoanInput")
$wsLoanInput.Range("B10").Select()

# --- Summary sheet: move selection from A6:XFD12 to A4 ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A4").Select()

# --- Repayment schedule sheet: insert a new column O (copy of N's formatting),
#     filled with 0 for rows 3-8, blank (formatted only) for row 2 ---
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Range("N2:N8").Copy()
$wsSchedule.Range("O2:O8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsSchedule.Range("O3:O8").Value = 0

# --- Edit Repayment Schedule sheet: move selection from A1:XFD1048576 to B7 ---
$wsEditSchedule = $wb.Worksheets.Item("Edit Repayment Schedule")
$wsEditSchedule.Range("B7").Select()

# Restore the originally active sheet/tab so only the intended selections change.
$originalActiveSheet.Activate()
